$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# RQ1 sheet: update Freecol (row 5) results, and small corrections to
# Ctags (row 3) and Jabref (row 7) rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("RQ1")

# Row 3 (Ctags): D3 13 -> 12, F3 2 -> 0, G3 13 -> 12 (E3/H3 are formulas, recalc automatically)
$ws1.Range("D3").Value = 12
$ws1.Range("F3").Value = 0
$ws1.Range("G3").Value = 12

# Row 5 (Freecol): previously blank C5/D5/F5/G5 (E5/H5 were #DIV/0!) -> now filled in
$ws1.Range("C5").Value = 1751
$ws1.Range("D5").Value = 93
$ws1.Range("F5").Value = 67
$ws1.Range("G5").Value = 93

# Row 7 (Jabref): D7 36 -> 35, F7 14 -> 8, G7 36 -> 35
$ws1.Range("D7").Value = 35
$ws1.Range("F7").Value = 8
$ws1.Range("G7").Value = 35

$ws1.Activate()
$ws1.Range("F5").Select() | Out-Null

# ---------------------------------------------------------------------------
# RQ2 sheet: fill in previously-blank B/C (and E/F for row 7) columns
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("RQ2")

# Row 3 (Ctags)
$ws2.Range("B3").Value = 32
$ws2.Range("C3").Value = 32
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 0

# Row 4 (Brlcad)
$ws2.Range("B4").Value = 33
$ws2.Range("C4").Value = 36
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 0

# Row 7 (Jabref)
$ws2.Range("B7").Value = 171
$ws2.Range("C7").Value = 198
$ws2.Range("E7").Value = 8
$ws2.Range("F7").Value = 9

$ws2.Activate()
$ws2.Range("F7").Select() | Out-Null

# ---------------------------------------------------------------------------
# RQ3 sheet: fill in previously-blank B/C/E/F columns
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("RQ3")

# Row 3 (Ctags)
$ws3.Range("B3").Value = 12
$ws3.Range("C3").Value = 12
$ws3.Range("E3").Value = 1
$ws3.Range("F3").Value = 1

# Row 4 (Brlcad)
$ws3.Range("B4").Value = 10
$ws3.Range("C4").Value = 10
$ws3.Range("E4").Value = 1
$ws3.Range("F4").Value = 1

# Row 7 (Jabref)
$ws3.Range("B7").Value = 35
$ws3.Range("C7").Value = 36
$ws3.Range("E7").Value = 3
$ws3.Range("F7").Value = 3

$ws3.Activate()
$ws3.Range("C7").Select() | Out-Null

# ---------------------------------------------------------------------------
# RQ4 sheet: fill in previously-blank B/C/E/F columns
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("RQ4")

# Row 3 (Ctags)
$ws4.Range("B3").Value = 12
$ws4.Range("C3").Value = 5
$ws4.Range("E3").Value = 1
$ws4.Range("F3").Value = 1

# Row 4 (Brlcad)
$ws4.Range("B4").Value = 10
$ws4.Range("C4").Value = 2
$ws4.Range("E4").Value = 1
$ws4.Range("F4").Value = 1

$ws4.Activate()
$ws4.Range("F5").Select() | Out-Null

# Leave RQ1 as the final active sheet/view (tabSelected="1" in RQ1 originally)
$ws1.Activate()
